$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("After")
$ws.Activate()
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 16.5
$ws.Range("F6").Value = 18.5
$ws.Range("G6").Value = 20
$ws.Range("I6").Value = 30
$ws.Range("H11").Select() | Out-Null
